# The header cell C1 previously read "Number of servings " (with a
# trailing space). The author replaced it with the shorter label
# "Servings".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "Servings"

# The author also left the selection/active cell on C19 when they
# saved the file (was A17 before).
$ws.Range("C19").Select()
